$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 348
$ws.Range("I4").Value = 348
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 348
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -234
$ws.Range("N4").ClearContents()

$ws.Range("H21").Value = 27509.5
$ws.Range("I21").Value = 20019
$ws.Range("J21").Value = 35000
$ws.Range("K21").Value = 20019
$ws.Range("L21").Value = 35000
$ws.Range("M21").Value = -19551
$ws.Range("N21").Value = -35936

$ws.Range("H23").Value = 27509.5
$ws.Range("I23").Value = 20019
$ws.Range("J23").Value = 35000
$ws.Range("K23").Value = 20019
$ws.Range("L23").Value = 35000
$ws.Range("M23").Value = -19785
$ws.Range("N23").Value = -35468

$ws.Range("H116").Value = 2763.6
$ws.Range("I116").Value = 2310.476
$ws.Range("J116").Value = 3443.2856
$ws.Range("K116").Value = 2310.476
$ws.Range("L116").Value = 3443.2856
$ws.Range("M116").Value = 1131.524
$ws.Range("N116").Value = -10327.2856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 23483.334
$ws.Range("I19").Value = 225
$ws.Range("K19").Value = 225
$ws.Range("M19").Value = 4

$ws.Range("H32").Value = 5581.565
$ws.Range("I32").Value = 4516.7744
$ws.Range("J32").Value = 10981.571
$ws.Range("K32").Value = 4516.7744
$ws.Range("L32").Value = 10981.571
$ws.Range("M32").Value = -4229.7744
$ws.Range("N32").Value = -11555.571

$ws.Range("H97").Value = 648.63635
$ws.Range("I97").Value = 640.8
$ws.Range("K97").Value = 640.8
$ws.Range("M97").Value = -144.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 770.85187
$ws.Range("I94").Value = 558.5
$ws.Range("J94").Value = 1195.5555
$ws.Range("K94").Value = 558.5
$ws.Range("L94").Value = 1195.5555
$ws.Range("M94").Value = -107.5
$ws.Range("N94").Value = -2097.5555

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 10002
$ws.Range("J4").Value = 10002
$ws.Range("L4").Value = 10002
$ws.Range("N4").Value = -10226

$ws.Range("H22").Value = 1565.3334
$ws.Range("I22").Value = 1881.6666
$ws.Range("J22").Value = 932.6667
$ws.Range("K22").Value = 1881.6666
$ws.Range("L22").Value = 932.6667
$ws.Range("M22").Value = -1531.6666
$ws.Range("N22").Value = -1632.6667

$ws.Range("H25").Value = 8337
$ws.Range("I25").Value = 5005.5
$ws.Range("J25").Value = 15000
$ws.Range("K25").Value = 5005.5
$ws.Range("L25").Value = 15000
$ws.Range("M25").Value = -4831.5
$ws.Range("N25").Value = -15348

$ws.Range("H58").Value = 1277.8235
$ws.Range("I58").Value = 1172.8462
$ws.Range("J58").Value = 1619
$ws.Range("K58").Value = 1172.8462
$ws.Range("L58").Value = 1619
$ws.Range("M58").Value = -969.8462
$ws.Range("N58").Value = -2025

$ws.Range("H134").Value = 2015
$ws.Range("I134").Value = 2156.9524
$ws.Range("J134").Value = 1269.75
$ws.Range("K134").Value = 6470.8572
$ws.Range("L134").Value = 3809.25
$ws.Range("M134").Value = -3935.8572
$ws.Range("N134").Value = -8879.25

$ws.Range("H136").Value = 1277.8235
$ws.Range("I136").Value = 1172.8462
$ws.Range("J136").Value = 1619
$ws.Range("K136").Value = 3518.5386
$ws.Range("L136").Value = 4857
$ws.Range("M136").Value = -968.5385999999999
$ws.Range("N136").Value = -9957

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 916.6667
$ws.Range("I25").Value = 275
$ws.Range("J25").Value = 2200
$ws.Range("K25").Value = 825
$ws.Range("L25").Value = 6600
$ws.Range("M25").Value = -656
$ws.Range("N25").Value = -6938

$ws.Range("H30").Value = 916.6667
$ws.Range("I30").Value = 275
$ws.Range("J30").Value = 2200
$ws.Range("K30").Value = 825
$ws.Range("L30").Value = 6600
$ws.Range("M30").Value = -723
$ws.Range("N30").Value = -6804

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 15831.429
$ws.Range("I17").Value = 120
$ws.Range("J17").Value = 18450
$ws.Range("K17").Value = 120
$ws.Range("L17").Value = 18450
$ws.Range("M17").Value = 48
$ws.Range("N17").Value = -18786

$ws.Range("H70").Value = 3927.2856
$ws.Range("I70").Value = 3878.85
$ws.Range("J70").Value = 4142.5557
$ws.Range("K70").Value = 3878.85
$ws.Range("L70").Value = 4142.5557
$ws.Range("M70").Value = -3608.85
$ws.Range("N70").Value = -4682.5557

$ws.Range("H73").Value = 3927.2856
$ws.Range("I73").Value = 3878.85
$ws.Range("J73").Value = 4142.5557
$ws.Range("K73").Value = 3878.85
$ws.Range("L73").Value = 4142.5557
$ws.Range("M73").Value = -2942.85
$ws.Range("N73").Value = -6014.5557

$ws.Range("H80").Value = 2512.3333
$ws.Range("I80").Value = 2451.25
$ws.Range("J80").Value = 2561.2
$ws.Range("K80").Value = 2451.25
$ws.Range("L80").Value = 2561.2
$ws.Range("M80").Value = -1453.25
$ws.Range("N80").Value = -4557.2

$ws.Range("H83").Value = 2512.3333
$ws.Range("I83").Value = 2451.25
$ws.Range("J83").Value = 2561.2
$ws.Range("K83").Value = 12256.25
$ws.Range("L83").Value = 12806
$ws.Range("M83").Value = -7264.25
$ws.Range("N83").Value = -22790

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 19500
$ws.Range("J2").Value = 19500
$ws.Range("L2").Value = 19500
$ws.Range("N2").Value = -19724

$ws.Range("H18").Value = 60000
$ws.Range("J18").Value = 60000
$ws.Range("L18").Value = 60000
$ws.Range("N18").Value = -60344

$ws.Range("H40").Value = 3585.5454
$ws.Range("I40").Value = 3122.0715
$ws.Range("J40").Value = 6181
$ws.Range("K40").Value = 3122.0715
$ws.Range("L40").Value = 6181
$ws.Range("M40").Value = -2986.0715
$ws.Range("N40").Value = -6453

$ws.Range("H61").Value = 1812.125
$ws.Range("I61").Value = 1529.4
$ws.Range("J61").Value = 2283.3333
$ws.Range("K61").Value = 1529.4
$ws.Range("L61").Value = 2283.3333
$ws.Range("M61").Value = -1327.4
$ws.Range("N61").Value = -2687.3333

$ws.Range("H113").Value = 1812.125
$ws.Range("I113").Value = 1529.4
$ws.Range("J113").Value = 2283.3333
$ws.Range("K113").Value = 1529.4
$ws.Range("L113").Value = 2283.3333
$ws.Range("M113").Value = 640.5999999999999
$ws.Range("N113").Value = -6623.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 172350
$ws.Range("I20").Value = 503050
$ws.Range("J20").Value = 7000
$ws.Range("K20").Value = 503050
$ws.Range("L20").Value = 7000
$ws.Range("M20").Value = -502810
$ws.Range("N20").Value = -7480

$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H132").Value = 455.48276
$ws.Range("I132").Value = 383.48148
$ws.Range("J132").Value = 1427.5
$ws.Range("K132").Value = 1150.44444
$ws.Range("L132").Value = 4282.5
$ws.Range("M132").Value = 1379.55556
$ws.Range("N132").Value = -9342.5
